$wb = $excel.ActiveWorkbook

# Remember the originally active sheet so we can restore it at the end
# (adding a new sheet below makes it the active one).
$origActive = $wb.ActiveSheet.Name

# --- 1. Insert the new "2022-Q1" sheet right before "总计" ---
$totalSheet = $wb.Sheets.Item("总计")
$newSheet = $wb.Worksheets.Add($totalSheet)
$newSheet.Name = "2022-Q1"

# Match the sheetPr / page setup used by the other quarterly sheets.
$newSheet.Outline.SummaryRow = 1
$newSheet.Outline.SummaryColumn = 1
$newSheet.PageSetup.LeftMargin = 54
$newSheet.PageSetup.RightMargin = 54
$newSheet.PageSetup.TopMargin = 72
$newSheet.PageSetup.BottomMargin = 72
$newSheet.PageSetup.HeaderMargin = 36
$newSheet.PageSetup.FooterMargin = 36

# Copy header-row / index-column formatting from an existing quarterly
# sheet (they all share the same style) so the new sheet matches exactly.
$template = $wb.Sheets.Item("2021-Q4")
$template.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$template.Range("A2").Copy()
$newSheet.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Header row text.
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Data row. B:G look numeric ("005126", "0.10", ...) but must stay text
# (matching every other quarterly sheet), so force text format before
# assigning, then drop back to the Normal style so no stray numbered
# format lingers on the cells.
$newSheet.Range("A2").Value = 0

$textRange = $newSheet.Range("B2:G2")
$textRange.NumberFormat = "@"
$newSheet.Range("B2").Value = "005126"
$newSheet.Range("C2").Value = "银河量化稳进混合"
$newSheet.Range("D2").Value = "0.10"
$newSheet.Range("E2").Value = "78.20"
$newSheet.Range("F2").Value = "2.28"
$newSheet.Range("G2").Value = "0.0023"
$textRange.Style = "Normal"

$newSheet.Range("H2").Value = 2

# --- 2. Insert a new first data row into "总计" sheet for 2022-Q1 ---
# Values are pushed down one row. Literal values are used (rather than
# reading-then-rewriting the existing cells) so the original numeric
# literals ("1.13", "1.37", "1.06") round-trip byte-for-byte instead of
# picking up binary floating-point noise from a COM Value2 round trip.
$total = $wb.Sheets.Item("总计")

$total.Range("B5").Value = "2021-Q2"
$total.Range("C5").Value = 5
$total.Range("D5").Value = 1.06

$total.Range("B4").Value = "2021-Q3"
$total.Range("C4").Value = 5
$total.Range("D4").Value = 1.37

$total.Range("B3").Value = "2021-Q4"
$total.Range("C3").Value = 8
$total.Range("D3").Value = 1.13

$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 1
$total.Range("D2").Value = 0

$total.Range("A2").Copy()
$total.Range("A5").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$total.Range("A5").Value = 3

# --- 3. Restore original active sheet/tab ---
$wb.Sheets.Item($origActive).Activate()
